$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 156 - "Journeymen Description" (Technology Description for the Journeymen's Guild)
$ws.Cells.Item(156, 2).Value = "Journeymen Description"
$ws.Cells.Item(156, 4).Value = "The Journeymen's Guild is a special organization of laborers that train craftsmen and masons. Having their guild in your kingdom reduces the costs of constructing new building components and provides highly skilled workers for the construction of new settlements. Journeymen can be recruited in any Mareten settlement."
$ws.Cells.Item(156, 7).Value = "Technology Description"

# Row 157 - "Journeyman" (ProperName for the PIONEER.INI unit)
$ws.Cells.Item(157, 2).Value = "ObjectData ProperName"
$ws.Cells.Item(157, 4).Value = "Journeyman"
$ws.Cells.Item(157, 7).Value = "PIONEER.INI ProperName"

# Row 158 - Journeymen unit description
$ws.Cells.Item(158, 2).Value = "UnitData Description"
$ws.Cells.Item(158, 4).Value = "Journeymen are expert craftsmen and builders. They can construct new settlements 50% faster than the less well-trained Settlers."
$ws.Cells.Item(158, 7).Value = "PIONEER.INI Description"

# Match the row heights to the re-wrapped text (as produced by the original author's save)
$ws.Rows.Item(156).RowHeight = 51.8
$ws.Rows.Item(157).RowHeight = 26.5
$ws.Rows.Item(158).RowHeight = 26.5

# Restore the view's scroll position / active selection
$win = $excel.ActiveWindow
$win.ScrollRow = 148
$win.ScrollColumn = 5
$ws.Range("H158").Select()
